$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseDate = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0

# Row 2 values (swap with row 4)
$ws.Range("D2").Value = $baseDate.AddDays(44273)
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 14000
$ws.Range("M2").Value = 14000
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 233

# Row 4 values (swap with row 2)
$ws.Range("D4").Value = $baseDate.AddDays(44291)
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("O4").Value = "Limache"
$ws.Range("P4").Value = 183
